# Fruta / hortaliza, semanal
#
# Inserts one new daily price observation row into the "Brócoli" sheet.
# The new record is inserted at row 173 (pushing the existing rows 173-240
# down to 174-241) so the sheet's chronological/insert order is preserved,
# matching the source data feed's weekly append behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shifts rows 173..240 down to 174..241,
# which also grows the used range from R240 to R241.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(173, 1).Value  = 11
$ws.Cells.Item(173, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(173, 3).Value  = "Bíobío"
$ws.Cells.Item(173, 4).Value  = 44636
$ws.Cells.Item(173, 5).Value  = 8
$ws.Cells.Item(173, 6).Value  = 100112023
$ws.Cells.Item(173, 7).Value  = "Brócoli"
$ws.Cells.Item(173, 8).Value  = "Sin especificar"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 150
$ws.Cells.Item(173, 11).Value = 800
$ws.Cells.Item(173, 12).Value = 1000
$ws.Cells.Item(173, 13).Value = 893
$ws.Cells.Item(173, 14).Value = "$/unidad"
$ws.Cells.Item(173, 15).Value = "Región Metropolitana"
$ws.Cells.Item(173, 16).Value = 893
$ws.Cells.Item(173, 17).Value = 1
$ws.Cells.Item(173, 18).Value = "Hortaliza"
